# Auto-generated script applying market-price refresh updates to the Jenova_Profits workbook
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-crafting profit sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 144186.14
$ws.Range("J70").Value = 144186.14
$ws.Range("L70").Value = 432558.42
$ws.Range("N70").Value = -433098.42
$ws.Range("H73").Value = 144186.14
$ws.Range("J73").Value = 144186.14
$ws.Range("L73").Value = 432558.42
$ws.Range("N73").Value = -434430.42
$ws.Range("H132").Value = 1567.6586
$ws.Range("I132").Value = 1406.85
$ws.Range("K132").Value = 4220.549999999999
$ws.Range("M132").Value = -1690.549999999999
$ws.Range("H138").Value = 4051.889
$ws.Range("J138").Value = 4229.5713
$ws.Range("L138").Value = 12688.7139
$ws.Range("N138").Value = -22968.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 100342.45
$ws.Range("I2").Value = 110346.6
$ws.Range("K2").Value = 110346.6
$ws.Range("M2").Value = -110233.6
$ws.Range("H32").Value = 3592.3418
$ws.Range("I32").Value = 3231.4133
$ws.Range("K32").Value = 3231.4133
$ws.Range("M32").Value = -2944.4133
$ws.Range("H97").Value = 938.6957
$ws.Range("I97").Value = 1019.4211
$ws.Range("J97").Value = 555.25
$ws.Range("K97").Value = 1019.4211
$ws.Range("L97").Value = 555.25
$ws.Range("M97").Value = -523.4211
$ws.Range("N97").Value = -1547.25
$ws.Range("H116").Value = 100342.45
$ws.Range("I116").Value = 110346.6
$ws.Range("K116").Value = 110346.6
$ws.Range("M116").Value = -108052.6
$ws.Range("H122").Value = 4741.6206
$ws.Range("I122").Value = 5409.25
$ws.Range("J122").Value = 3919.923
$ws.Range("K122").Value = 16227.75
$ws.Range("L122").Value = 11759.769
$ws.Range("M122").Value = -13777.75
$ws.Range("N122").Value = -16659.769
$ws.Range("H132").Value = 3310.6667
$ws.Range("I132").Value = 3313.8286
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 9941.485799999999
$ws.Range("L132").Value = 9600
$ws.Range("M132").Value = -7411.485799999999
$ws.Range("N132").Value = -14660

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 100342.45
$ws.Range("I3").Value = 110346.6
$ws.Range("K3").Value = 110346.6
$ws.Range("M3").Value = -110232.6
$ws.Range("H20").Value = 38464550
$ws.Range("I20").Value = 62502710
$ws.Range("J20").Value = 3495
$ws.Range("K20").Value = 62502710
$ws.Range("L20").Value = 3495
$ws.Range("M20").Value = -62502463
$ws.Range("N20").Value = -3989
$ws.Range("H86").Value = 740778.1
$ws.Range("I86").Value = 945794.3
$ws.Range("K86").Value = 945794.3
$ws.Range("M86").Value = -944671.3
$ws.Range("H89").Value = 740778.1
$ws.Range("I89").Value = 945794.3
$ws.Range("K89").Value = 4728971.5
$ws.Range("M89").Value = -4723355.5
$ws.Range("H94").Value = 898.5599999999999
$ws.Range("I94").Value = 1058.4
$ws.Range("J94").Value = 259.2
$ws.Range("K94").Value = 1058.4
$ws.Range("L94").Value = 259.2
$ws.Range("M94").Value = -607.4000000000001
$ws.Range("N94").Value = -1161.2
$ws.Range("H105").Value = 92198.37
$ws.Range("I105").Value = 92198.37
$ws.Range("K105").Value = 92198.37
$ws.Range("M105").Value = -90451.37
$ws.Range("H107").Value = 501763.1
$ws.Range("I107").Value = 1236
$ws.Range("J107").Value = 2003344.4
$ws.Range("K107").Value = 1236
$ws.Range("L107").Value = 2003344.4
$ws.Range("M107").Value = 684
$ws.Range("N107").Value = -2007184.4
$ws.Range("H134").Value = 44645.332
$ws.Range("I134").Value = 3108.1738
$ws.Range("K134").Value = 9324.5214
$ws.Range("M134").Value = -6789.5214

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 315.625
$ws.Range("I2").Value = 59.333332
$ws.Range("J2").Value = 469.4
$ws.Range("K2").Value = 355.999992
$ws.Range("L2").Value = 2816.4
$ws.Range("M2").Value = -242.999992
$ws.Range("N2").Value = -3042.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 142879550
$ws.Range("I70").Value = 19223.25
$ws.Range("K70").Value = 19223.25
$ws.Range("M70").Value = -18953.25
$ws.Range("H73").Value = 142879550
$ws.Range("I73").Value = 19223.25
$ws.Range("K73").Value = 19223.25
$ws.Range("M73").Value = -18287.25
$ws.Range("H80").Value = 2226538
$ws.Range("I80").Value = 1433705.9
$ws.Range("J80").Value = 5001450
$ws.Range("K80").Value = 1433705.9
$ws.Range("L80").Value = 5001450
$ws.Range("M80").Value = -1432707.9
$ws.Range("N80").Value = -5003446
$ws.Range("H83").Value = 2226538
$ws.Range("I83").Value = 1433705.9
$ws.Range("J83").Value = 5001450
$ws.Range("K83").Value = 7168529.5
$ws.Range("L83").Value = 25007250
$ws.Range("M83").Value = -7163537.5
$ws.Range("N83").Value = -25017234
$ws.Range("H102").Value = 2006.25
$ws.Range("J102").Value = 4333.3335
$ws.Range("L102").Value = 4333.3335
$ws.Range("N102").Value = -7577.3335
$ws.Range("H132").Value = 112250.2
$ws.Range("I132").Value = 15001.4
$ws.Range("J132").Value = 209499
$ws.Range("K132").Value = 45004.2
$ws.Range("L132").Value = 628497
$ws.Range("M132").Value = -42474.2
$ws.Range("N132").Value = -633557
$ws.Range("H135").Value = 333458340
$ws.Range("J135").Value = 333458340
$ws.Range("L135").Value = 333458340
$ws.Range("N135").Value = -333468480

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6665.524
$ws.Range("I7").Value = 6818.4
$ws.Range("K7").Value = 6818.4
$ws.Range("M7").Value = -6706.4
$ws.Range("H122").Value = 3015.1428
$ws.Range("I122").Value = 2825.2856
$ws.Range("J122").Value = 3584.7144
$ws.Range("K122").Value = 8475.856800000001
$ws.Range("L122").Value = 10754.1432
$ws.Range("M122").Value = -6025.856800000001
$ws.Range("N122").Value = -15654.1432
$ws.Range("H126").Value = 6665.524
$ws.Range("I126").Value = 6818.4
$ws.Range("K126").Value = 20455.2
$ws.Range("M126").Value = -17985.2
$ws.Range("H132").Value = 7197.4
$ws.Range("I132").Value = 6664.222
$ws.Range("K132").Value = 19992.666
$ws.Range("M132").Value = -17462.666
$ws.Range("H136").Value = 2006182.8
$ws.Range("I136").Value = 3337302.8
$ws.Range("K136").Value = 10011908.4
$ws.Range("M136").Value = -10009358.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H56").Value = 31061.25
$ws.Range("J56").Value = 31061.25
$ws.Range("L56").Value = 31061.25
$ws.Range("N56").Value = -32489.25
$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 20000
$ws.Range("K61").Value = 20000
$ws.Range("M61").Value = -19708
$ws.Range("H96").Value = 252749.75
$ws.Range("I96").Value = 335666.34
$ws.Range("J96").Value = 4000
$ws.Range("K96").Value = 335666.34
$ws.Range("L96").Value = 4000
$ws.Range("M96").Value = -334293.34
$ws.Range("N96").Value = -6746
$ws.Range("H122").Value = 25001226
$ws.Range("I122").Value = 29413032
$ws.Range("K122").Value = 88239096
$ws.Range("M122").Value = -88236646
$ws.Range("H132").Value = 40783.742
$ws.Range("I132").Value = 2357.158
$ws.Range("K132").Value = 7071.474
$ws.Range("M132").Value = -4541.474
$ws.Range("H136").Value = 10819584
$ws.Range("I136").Value = 11907147
$ws.Range("J136").Value = 668999.3
$ws.Range("K136").Value = 35721441
$ws.Range("L136").Value = 2006997.9
$ws.Range("M136").Value = -35718891
$ws.Range("N136").Value = -2012097.9
